# "fixed no alias bug": column I ("alias") was left holding a placeholder
# "?" for every data row. The real alias is just the row's ID (column A),
# which Excel shows as a merged/blank-filled block - i.e. each run of blank
# cells below an ID belongs to that same ID. Re-create that fill-down into
# column I, replacing every "?" with the correct alias.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$lastAlias = $null

for ($row = 2; $row -le $lastRow; $row++) {
    $idText = $ws.Cells.Item($row, 1).Text
    if ($idText -ne "") {
        $lastAlias = $idText
    }

    $cell = $ws.Cells.Item($row, 9)
    # Force text storage (matches column A's own storage, and keeps "14_1"
    # etc. intact) instead of letting Excel auto-coerce numeric-looking
    # aliases like "10" into a Number, then drop back to the default style
    # so no stray number-format is left behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $lastAlias
    $cell.Style = "Normal"
}
